$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect B2:E51 from Excel's automatic number/date coercion while we
# write literal text values (prices like "243.02" or "1.000" would
# otherwise be parsed as numbers).
$protect = $ws.Range("B2:E51")
$protect.NumberFormat = "@"

# Rows 2-22: update Price (D) and Volume(1h) (E) only
$ws.Range("D2").Value = "29.825.30"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").Value = "1.894.72"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "0.7976"
$ws.Range("E5").Value = "  -2.75%  "
$ws.Range("D6").Value = "243.02"
$ws.Range("E6").Value = "  +0.69%  "
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "0.3170"
$ws.Range("E8").Value = "  -1.33%  "
$ws.Range("E9").Value = "  -3.12%  "
$ws.Range("D10").Value = "0.07049"
$ws.Range("E10").Value = "  +0.59%  "
$ws.Range("D11").Value = "0.08080"
$ws.Range("E11").Value = "  +0.66%  "
$ws.Range("D12").Value = "0.7717"
$ws.Range("E12").Value = "  +3.73%  "
$ws.Range("D13").Value = "1.887.33"
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("D14").Value = "5.357"
$ws.Range("E14").Value = "  +3.33%  "
$ws.Range("D15").Value = "92.71"
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("D16").Value = "29.843.63"
$ws.Range("E16").Value = "  -0.21%  "
$ws.Range("D17").Value = "5.999"
$ws.Range("E17").Value = "  +2.17%  "
$ws.Range("D18").Value = "13.89"
$ws.Range("E18").Value = "  -0.55%  "
$ws.Range("D19").Value = "245.09"
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("D20").Value = "0.000007717"
$ws.Range("E20").Value = "  -0.33%  "
$ws.Range("D21").Value = "8.318"
$ws.Range("E21").Value = "  +20.58%  "
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.12%  "

# Rows 23-51: Coin (B), Link (C), Price (D), Volume(1h) (E) all shift/update
$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("B24").Value = "Stellar"
$ws.Range("C24").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D24").Value = "0.1649"
$ws.Range("E24").Value = "  +6.30%  "
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").Value = "9.349"
$ws.Range("E25").Value = "  +2.03%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "166.04"
$ws.Range("E26").Value = "  +0.16%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "18.73"
$ws.Range("E27").Value = "  -0.32%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "2.062"
$ws.Range("E28").Value = "  -0.62%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "1.401"
$ws.Range("E29").Value = "  +2.49%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "1.543"
$ws.Range("E30").Value = "  +1.47%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "4.460"
$ws.Range("E31").Value = "  +4.64%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "0.05694"
$ws.Range("E32").Value = "  +1.29%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "4.048"
$ws.Range("E33").Value = "  -0.35%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "1.265"
$ws.Range("E34").Value = "  -0.14%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "0.7396"
$ws.Range("E35").Value = "  +1.58%  "
$ws.Range("B36").Value = "Frax"
$ws.Range("C36").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D36").Value = "0.9989"
$ws.Range("E36").Value = "  -0.17%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "2.641"
$ws.Range("E37").Value = "  -3.02%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.01911"
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "2.784"
$ws.Range("E39").Value = "  +0.43%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "0.4417"
$ws.Range("E40").Value = "  +0.25%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "72.55"
$ws.Range("E41").Value = "  +1.29%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "5.824"
$ws.Range("E42").Value = "  -2.15%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "0.8412"
$ws.Range("E43").Value = "  -0.33%  "
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "1.033.16"
$ws.Range("E45").Value = "  +4.70%  "
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "103.25"
$ws.Range("E46").Value = "  +2.48%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "10.04"
$ws.Range("E47").Value = "  +3.91%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "1.875"
$ws.Range("E48").Value = "  +0.30%  "
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").Value = "7.434"
$ws.Range("E49").Value = "  -1.62%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.053.67"
$ws.Range("E50").Value = "  +0.58%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "1.532"
$ws.Range("E51").Value = "  +4.70%  "

# Remove the temporary Text number-format so cells return to the default style
$protect.ClearFormats()
